# Auto-generated Excel market-data refresh edits
# Applies the value updates captured in the commit diff for Sheets/Ragnarok_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 71434330
$ws.Range("J40").Value = 83339810
$ws.Range("L40").Value = 83339810
$ws.Range("N40").Value = -83340160

$ws.Range("H101").Value = 788.35
$ws.Range("I101").Value = 907.3333
$ws.Range("J101").Value = 691
$ws.Range("K101").Value = 2721.9999
$ws.Range("L101").Value = 2073
$ws.Range("M101").Value = -1099.9999
$ws.Range("N101").Value = -5317

$ws.Range("H107").Value = 1371.5555
$ws.Range("I107").Value = 501.46155
$ws.Range("K107").Value = 501.46155
$ws.Range("M107").Value = 1418.53845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8653.522000000001
$ws.Range("I32").Value = 8170.984
$ws.Range("K32").Value = 8170.984
$ws.Range("M32").Value = -7883.984

$ws.Range("H45").Value = 3816
$ws.Range("I45").Value = 2576
$ws.Range("J45").Value = 6296
$ws.Range("K45").Value = 2576
$ws.Range("L45").Value = 6296
$ws.Range("M45").Value = -2199
$ws.Range("N45").Value = -7050

$ws.Range("H61").Value = 5839012.5
$ws.Range("I61").Value = 6672621.5
$ws.Range("K61").Value = 6672621.5
$ws.Range("M61").Value = -6672409.5

$ws.Range("H74").Value = 1655.6897
$ws.Range("I74").Value = 1563.5186
$ws.Range("K74").Value = 1563.5186
$ws.Range("M74").Value = -689.5186000000001

$ws.Range("H77").Value = 1655.6897
$ws.Range("I77").Value = 1563.5186
$ws.Range("K77").Value = 7817.593000000001
$ws.Range("M77").Value = -3449.593000000001

$ws.Range("H122").Value = 2086.75
$ws.Range("I122").Value = 1824.6
$ws.Range("K122").Value = 5473.799999999999
$ws.Range("M122").Value = -3023.799999999999

$ws.Range("H132").Value = 6670078.5
$ws.Range("I132").Value = 3475
$ws.Range("J132").Value = 33336494
$ws.Range("K132").Value = 10425
$ws.Range("L132").Value = 100009482
$ws.Range("M132").Value = -7895
$ws.Range("N132").Value = -100014542

$ws.Range("H136").Value = 5839012.5
$ws.Range("I136").Value = 6672621.5
$ws.Range("K136").Value = 20017864.5
$ws.Range("M136").Value = -20015314.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4860.9
$ws.Range("I94").Value = 4860.9
$ws.Range("K94").Value = 4860.9
$ws.Range("M94").Value = -4409.9

$ws.Range("H122").Value = 46065.848
$ws.Range("J122").Value = 45693.75
$ws.Range("L122").Value = 45693.75
$ws.Range("N122").Value = -55493.75

$ws.Range("H127").Value = 54900
$ws.Range("J127").Value = 54900
$ws.Range("L127").Value = 54900
$ws.Range("N127").Value = -64820

$ws.Range("H134").Value = 7693248
$ws.Range("I134").Value = 684.2727
$ws.Range("J134").Value = 50002350
$ws.Range("K134").Value = 2052.8181
$ws.Range("L134").Value = 150007050
$ws.Range("M134").Value = 482.1819
$ws.Range("N134").Value = -150012120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5883363.5
$ws.Range("I16").Value = 9091671
$ws.Range("J16").Value = 1465.6666
$ws.Range("K16").Value = 9091671
$ws.Range("L16").Value = 1465.6666
$ws.Range("M16").Value = -9091384
$ws.Range("N16").Value = -2039.6666

$ws.Range("H105").Value = 6114.6
$ws.Range("I105").Value = 1171
$ws.Range("K105").Value = 1171
$ws.Range("M105").Value = 576

$ws.Range("H113").Value = 5883363.5
$ws.Range("I113").Value = 9091671
$ws.Range("J113").Value = 1465.6666
$ws.Range("K113").Value = 9091671
$ws.Range("L113").Value = 1465.6666
$ws.Range("M113").Value = -9089501
$ws.Range("N113").Value = -5805.6666

$ws.Range("H122").Value = 1856.1034
$ws.Range("I122").Value = 1219.4546
$ws.Range("J122").Value = 3857
$ws.Range("K122").Value = 3658.3638
$ws.Range("L122").Value = 11571
$ws.Range("M122").Value = -1208.3638
$ws.Range("N122").Value = -16471

$ws.Range("H123").Value = 62499.668
$ws.Range("J123").Value = 62499.668
$ws.Range("L123").Value = 62499.668
$ws.Range("N123").Value = -72299.66800000001

$ws.Range("H132").Value = 3300.7778
$ws.Range("I132").Value = 3400.875
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 10202.625
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -7672.625
$ws.Range("N132").Value = -12560

$ws.Range("H134").Value = 3478
$ws.Range("I134").Value = 2890
$ws.Range("J134").Value = 4507
$ws.Range("K134").Value = 8670
$ws.Range("L134").Value = 13521
$ws.Range("M134").Value = -6135
$ws.Range("N134").Value = -18591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 138.42857
$ws.Range("I2").Value = 117.5
$ws.Range("J2").Value = 166.33333
$ws.Range("K2").Value = 705
$ws.Range("L2").Value = 997.9999799999999
$ws.Range("M2").Value = -592
$ws.Range("N2").Value = -1223.99998

$ws.Range("H26").Value = 4336.8887
$ws.Range("I26").Value = 919.8333
$ws.Range("K26").Value = 2759.4999
$ws.Range("M26").Value = -2471.4999

$ws.Range("H33").Value = 7699720
$ws.Range("J33").Value = 17324172
$ws.Range("L33").Value = 103945032
$ws.Range("N33").Value = -103945598

$ws.Range("H44").Value = 22166
$ws.Range("J44").Value = 22166
$ws.Range("L44").Value = 66498
$ws.Range("N44").Value = -67294

$ws.Range("H87").Value = 39999.332
$ws.Range("I87").Value = 9999.5
$ws.Range("K87").Value = 29998.5
$ws.Range("M87").Value = -28750.5

$ws.Range("H90").Value = 39999.332
$ws.Range("I90").Value = 9999.5
$ws.Range("K90").Value = 89995.5
$ws.Range("M90").Value = -83755.5

$ws.Range("H137").Value = 9366.066000000001
$ws.Range("I137").Value = 5082.5
$ws.Range("J137").Value = 10923.728
$ws.Range("K137").Value = 15247.5
$ws.Range("L137").Value = 32771.18399999999
$ws.Range("M137").Value = -10147.5
$ws.Range("N137").Value = -42971.18399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1841909.8
$ws.Range("I122").Value = 3005942.8
$ws.Range("J122").Value = 12715.143
$ws.Range("K122").Value = 9017828.399999999
$ws.Range("L122").Value = 38145.429
$ws.Range("M122").Value = -9015378.399999999
$ws.Range("N122").Value = -43045.429

$ws.Range("H132").Value = 12506024
$ws.Range("I132").Value = 7237.8
$ws.Range("J132").Value = 33337332
$ws.Range("K132").Value = 21713.4
$ws.Range("L132").Value = 100011996
$ws.Range("M132").Value = -19183.4
$ws.Range("N132").Value = -100017056

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8090.129
$ws.Range("I7").Value = 8255.944
$ws.Range("J7").Value = 7860.5386
$ws.Range("K7").Value = 8255.944
$ws.Range("L7").Value = 7860.5386
$ws.Range("M7").Value = -8143.944
$ws.Range("N7").Value = -8084.5386

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H46").Value = 785.73334
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376

$ws.Range("H126").Value = 8090.129
$ws.Range("I126").Value = 8255.944
$ws.Range("J126").Value = 7860.5386
$ws.Range("K126").Value = 24767.832
$ws.Range("L126").Value = 23581.6158
$ws.Range("M126").Value = -22297.832
$ws.Range("N126").Value = -28521.6158

$ws.Range("H128").Value = 68330.336
$ws.Range("J128").Value = 68330.336
$ws.Range("L128").Value = 68330.336
$ws.Range("N128").Value = -78290.336

$ws.Range("H131").Value = 99499.5
$ws.Range("J131").Value = 99499.5
$ws.Range("L131").Value = 99499.5
$ws.Range("N131").Value = -109579.5

$ws.Range("H132").Value = 4603.654
$ws.Range("I132").Value = 2317.3076
$ws.Range("J132").Value = 6890
$ws.Range("K132").Value = 6951.9228
$ws.Range("L132").Value = 20670
$ws.Range("M132").Value = -4421.9228
$ws.Range("N132").Value = -25730

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 5399.8335
$ws.Range("I23").Value = 299.66666
$ws.Range("J23").Value = 10500
$ws.Range("K23").Value = 299.66666
$ws.Range("L23").Value = 10500
$ws.Range("M23").Value = -70.66665999999998
$ws.Range("N23").Value = -10958

$ws.Range("H107").Value = 2821.7144
$ws.Range("I107").Value = 1857.9524
$ws.Range("J107").Value = 5713
$ws.Range("K107").Value = 5573.857199999999
$ws.Range("L107").Value = 17139
$ws.Range("M107").Value = -3653.857199999999
$ws.Range("N107").Value = -20979

$ws.Range("H122").Value = 1654.4138
$ws.Range("I122").Value = 1307
$ws.Range("J122").Value = 4665.3335
$ws.Range("K122").Value = 3921
$ws.Range("L122").Value = 13996.0005
$ws.Range("M122").Value = -1471
$ws.Range("N122").Value = -18896.0005

$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -109799

$ws.Range("H126").Value = 6896.875
$ws.Range("I126").Value = 7398.5713
$ws.Range("K126").Value = 22195.7139
$ws.Range("M126").Value = -19725.7139

$ws.Range("H136").Value = 180263.39
$ws.Range("I136").Value = 1672.0682
$ws.Range("K136").Value = 5016.2046
$ws.Range("M136").Value = -2466.2046

$ws.Range("H140").Value = 58789.4
$ws.Range("J140").Value = 58789.4
$ws.Range("L140").Value = 58789.4
$ws.Range("N140").Value = -69149.39999999999

$ws.Range("H141").Value = 93883
$ws.Range("J141").Value = 93883
$ws.Range("L141").Value = 93883
$ws.Range("N141").Value = -104243
